$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 62 (pushes existing rows 62..135 down to 63..136,
# Excel auto-extends the used range / dimension to R136).
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row 62 with the new weekly record.
$ws.Range("A62").Value = 9
$ws.Range("B62").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C62").Value = "Metropolitana"
$ws.Range("D62").Value = 45117
$ws.Range("E62").Value = 13
$ws.Range("F62").Value = 100114007
$ws.Range("G62").Value = "Jengibre"
$ws.Range("H62").Value = "Sin especificar"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 520
$ws.Range("K62").Value = 18000
$ws.Range("L62").Value = 20000
$ws.Range("M62").Value = 19000
$ws.Range("N62").Value = "$/caja 13 kilos"
$ws.Range("O62").Value = "Perú"
$ws.Range("P62").Value = 1462
$ws.Range("Q62").Value = 13
$ws.Range("R62").Value = "Hortaliza"
